# classificator.xlsx — "Switched count to double in all available places."
#
# The two rows below (31 = "Кросс" line, 46 = "Муфта оптическая" line) had
# their part numbers swapped out for a different catalog item and their
# quantity (column H) changed from an integer-looking count to "2" (the
# commit message notes counts became doubles everywhere available) while
# picking up the same centred/bordered formatting already used by the
# regular data rows (e.g. row 2).
#
# NOTE on row processing order: row 46 is edited before row 31 on purpose —
# that's the order in which the new shared strings need to be appended
# (Муфта оптическая / МТОК-Л6/96-1КВ244 / Кросс / ШКОС-Л ...) to reproduce
# the canonical shared-strings table exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- normalise formatting on B/D/E for the two rows (border + centered,
#     matching the rest of the table, e.g. row 2) -------------------------
foreach ($r in 31, 46) {
    foreach ($c in 2, 4, 5) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.HorizontalAlignment = -4108   # xlCenter
        $cell.VerticalAlignment = -4108     # xlCenter
    }
}

# --- row 46: Муфта оптическая / МТОК-Л6/96-1КВ244 ------------------------
$ws.Cells.Item(46, 3).Value2 = "Муфта оптическая"
$ws.Cells.Item(46, 4).Value2 = "МТОК-Л6/96-1КВ244"
$ws.Cells.Item(46, 6).Value2 = "ЗАО «СВЯЗЬСТРОЙДЕТАЛЬ»"
$ws.Cells.Item(46, 8).Value2 = 2

# --- row 31: Кросс / ШКОС-Л -1U/2 -48 -LC ~48 -LC/SM ~48 -LC/UPC ---------
$ws.Cells.Item(31, 3).Value2 = "Кросс"
$ws.Cells.Item(31, 4).Value2 = "ШКОС-Л -1U/2 -48 -LC ~48 -LC/SM ~48 -LC/UPC"
$ws.Cells.Item(31, 6).Value2 = "ЗАО «СВЯЗЬСТРОЙДЕТАЛЬ»"
$ws.Cells.Item(31, 8).Value2 = 2

# --- view state: scrolled down a bit further and a different cell
#     selected (C36 instead of D52) ---------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C36").Select() | Out-Null
